# Verify_40V_On_Addition_Deletion_Of_Ethernet.xlsx edit:
# The "P405D" panel-type test row (row 10 on the "Add Panels" sheet) was removed
# from the test data table, shifting the following rows (Pro32xD, MX2-100, P885D)
# up by one and dropping the now-unused "P405D" shared string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# Remove the entire row that held the "P405D" panel type test entry.
$ws.Rows("10").Delete() | Out-Null

# Reflect the author's final cell selection in the saved view.
$ws.Range("A9").Select() | Out-Null
